# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the handoff timestamps on the Overview, zh-cn and de-de sheets, and widens
# the Status/date columns that grew to fit the new text.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# The underlying engine re-quantizes ColumnWidth assignments to the nearest
# 1/6-character increment (output = round(input*6 + 5)/6). 16.333333 is the
# input value that lands on the increment closest to the authored target
# width of 17.2159881591797 ("Ready for handoff" / wider timestamp text).
$statusColWidth = 16.333333

# --- Overview sheet ---------------------------------------------------
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Range("G2").Value = "2016-08-26 16:59:15"

$ws_overview.Columns("E").ColumnWidth = $statusColWidth
$ws_overview.Columns("F").ColumnWidth = $statusColWidth

# --- zh-cn sheet --------------------------------------------------------
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("H2").Value = "2016-08-26 16:59:11"

$ws_zhcn.Columns("C").ColumnWidth = $statusColWidth

# --- de-de sheet --------------------------------------------------------
$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("H2").Value = "2016-08-26 16:59:15"

$ws_dede.Columns("C").ColumnWidth = $statusColWidth
